$d = $word.ActiveDocument

$d.Content.Find.Execute("[1] 4.817360 5.385498", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "[1] 4.738957 5.455328", 2)

$d.Content.Find.Execute("[1] 0.954", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "[1] 0.952", 2)

$d.Content.Find.Execute("[1] 0.885", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "[1] 0.872", 2)

$d.Content.Find.Execute("[1] 4.879784 5.742216", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "[1] 4.494485 4.953515", 2)
